# Disambiguate the four previously-identical {blood_test_results}
# placeholders into distinct, test-specific placeholders so the
# template no longer raises an UndefinedError / mixes up results.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Human immunodeficiency virus (HIV), types 1 and 2: {blood_test_results}",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Human immunodeficiency virus (HIV), types 1 and 2: {hiv_results}",
    2)

$d.Content.Find.Execute(
    "Hepatitis B virus (HBV): {blood_test_results}",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Hepatitis B virus (HBV): {hbv_results}",
    2)

$d.Content.Find.Execute(
    "Hepatitis C virus (HCV): {blood_test_results}",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Hepatitis C virus (HCV): {hcv_results}",
    2)

$d.Content.Find.Execute(
    "Treponema pallidum (syphilis) through VDRL: {blood_test_results}",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Treponema pallidum (syphilis) through VDRL: {vdrl_results}",
    2)
